# "Add Doc Type Test Document to Template"
#
# The "Types" sheet holds the list of valid Document Types (used as the
# source list for the Type column's data-validation dropdown on the
# "Documents" sheet). Add a new entry, "Test Document", positioned right
# after "Interface Specification" and before "End-User Documentation" -
# i.e. insert a new row 10 and push the remaining rows down by one.

$wb = $excel.ActiveWorkbook

$types = $wb.Worksheets.Item("Types")
$types.Activate()

# Insert a new row at position 10 (shifts rows 10-12 down to 11-13) and
# fill it in with the new document type.
$types.Rows("10:10").Insert()
$types.Cells.Item(10, 1).Value = "Test Document"

# Leave the cursor where the edit happened on the Types sheet...
$types.Range("E10").Select()

# ...then hop back to the Documents sheet (the sheet the workbook is
# normally viewed on) and leave the selection parked on C20, matching
# where the author's cursor ended up.
$documents = $wb.Worksheets.Item("Documents")
$documents.Activate()
$documents.Range("C20").Select()
